$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price (column D) values are plain decimal numbers. Excel would
# auto-convert such strings to numeric cells on assignment, but the source
# data must remain text (t="inlineStr"/shared string), matching the rest of
# the column. Temporarily mark those cells as Text so the assignment keeps
# them as strings, then restore the original (default) cell formatting.
$textCells = @("D5", "D6", "D8", "D9", "D11", "D13", "D14", "D19", "D20", "D21", "D22", "D26", "D27", "D28", "D33", "D35", "D36", "D37", "D38", "D39", "D42", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range('D2').Value = '91.556.62'
$ws.Range('E2').Value = '  +0.99%  '
# Row 3
$ws.Range('D3').Value = '3.128.28'
$ws.Range('E3').Value = '  +0.61%  '
# Row 4
$ws.Range('E4').Value = '  -0.01%  '
# Row 5
$ws.Range('D5').Value = '241.57'
$ws.Range('E5').Value = '  -0.47%  '
# Row 6
$ws.Range('D6').Value = '618.55'
$ws.Range('E6').Value = '  -0.81%  '
# Row 7
$ws.Range('E7').Value = '  -5.67%  '
# Row 8
$ws.Range('D8').Value = '0.390'
$ws.Range('E8').Value = '  +5.26%  '
# Row 9
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.04%  '
# Row 10
$ws.Range('D10').Value = '3.128.91'
$ws.Range('E10').Value = '  +0.72%  '
# Row 11
$ws.Range('D11').Value = '0.751'
$ws.Range('E11').Value = '  -1.00%  '
# Row 12
$ws.Range('E12').Value = '  +0.30%  '
# Row 13
$ws.Range('D13').Value = '0.0000255'
$ws.Range('E13').Value = '  +1.27%  '
# Row 14
$ws.Range('D14').Value = '35.18'
$ws.Range('E14').Value = '  -0.36%  '
# Row 15
$ws.Range('E15').Value = '  +2.22%  '
# Row 16
$ws.Range('D16').Value = '91.392.97'
$ws.Range('E16').Value = '  +0.97%  '
# Row 17
$ws.Range('E17').Value = '  +0.82%  '
# Row 18
$ws.Range('D18').Value = '3.161.79'
$ws.Range('E18').Value = '  +2.36%  '
# Row 19
$ws.Range('D19').Value = '3.77'
$ws.Range('E19').Value = '  -0.39%  '
# Row 20
$ws.Range('D20').Value = '14.97'
$ws.Range('E20').Value = '  +4.34%  '
# Row 21
$ws.Range('D21').Value = '5.92'
$ws.Range('E21').Value = '  +1.87%  '
# Row 22
$ws.Range('D22').Value = '456.87'
$ws.Range('E22').Value = '  +1.84%  '
# Row 24
$ws.Range('E24').Value = '  +1.69%  '
# Row 25
$ws.Range('E25').Value = '  +0.85%  '
# Row 26
$ws.Range('D26').Value = '88.96'
$ws.Range('E26').Value = '  -4.99%  '
# Row 27
$ws.Range('D27').Value = '11.79'
$ws.Range('E27').Value = '  -0.98%  '
# Row 28
$ws.Range('D28').Value = '0.152'
$ws.Range('E28').Value = '  +35.09%  '
# Row 29
$ws.Range('D29').Value = '3.312.28'
$ws.Range('E29').Value = '  +1.53%  '
# Row 30
$ws.Range('E30').Value = '  -0.06%  '
# Row 31
$ws.Range('E31').Value = '  +1.58%  '
# Row 32
$ws.Range('E32').Value = '  -5.14%  '
# Row 33
$ws.Range('D33').Value = '9.36'
$ws.Range('E33').Value = '  +2.84%  '
# Row 34
$ws.Range('E34').Value = '  +9.78%  '
# Row 35
$ws.Range('D35').Value = '26.42'
$ws.Range('E35').Value = '  -0.79%  '
# Row 36
$ws.Range('D36').Value = '7.46'
$ws.Range('E36').Value = '  -2.32%  '
# Row 37
$ws.Range('D37').Value = '1.96'
$ws.Range('E37').Value = '  +2.24%  '
# Row 38
$ws.Range('D38').Value = '3.95'
$ws.Range('E38').Value = '  -5.54%  '
# Row 39
$ws.Range('D39').Value = '492.52'
$ws.Range('E39').Value = '  -0.21%  '
# Row 40
$ws.Range('E40').Value = '  +1.98%  '
# Row 41
$ws.Range('E41').Value = '  +5.19%  '
# Row 42
$ws.Range('D42').Value = '3.40'
$ws.Range('E42').Value = '  -5.60%  '
# Row 43
$ws.Range('E43').Value = '  +0.22%  '
# Row 45
$ws.Range('E45').Value = '  -28.49%  '
# Row 46
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').Value = '1.93'
$ws.Range('E46').Value = '  +1.40%  '
# Row 47
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Value = '0.708'
$ws.Range('E47').Value = '  +3.19%  '
# Row 48
$ws.Range('D48').Value = '156.39'
$ws.Range('E48').Value = '  -0.62%  '
# Row 49
$ws.Range('E49').Value = '  +1.58%  '
# Row 50
$ws.Range('D50').Value = '4.46'
$ws.Range('E50').Value = '  -2.60%  '
# Row 51
$ws.Range('D51').Value = '0.0328'
$ws.Range('E51').Value = '  +5.25%  '

# Restore original (default, non-"@") cell formatting on the cells we
# temporarily marked as Text, by pasting the format from the adjacent
# Coin column (same row), which still carries the untouched default style.
foreach ($addr in $textCells) {
    $row = $addr.Substring(1)
    $ws.Range("B$row").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
